$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Strip out the bank-account-specific values (they become blank template
#    cells). The labels in column A stay; the shared-string pool and the
#    value indices naturally re-pack once the unused strings are gone.
# ---------------------------------------------------------------------------
$ws.Range("B18").Value = $null   # was: AO <<BANK>>
$ws.Range("B19").Value = $null   # was: 0000000000000000000
$ws.Range("B20").Value = $null   # was: Russia
$ws.Range("B21").Value = $null   # was: 1 Lenina str., Moscow, ...
$ws.Range("B28").Value = $null   # was: IP Ivanov Ivan Ivanovich
$ws.Range("B29").Value = $null   # was: PR. LENINA, D. 1, KV. 1 ... (rich text)

# ---------------------------------------------------------------------------
# 2) B28:D28 gets a new centered style (horizontal=center, vertical=center)
#    instead of the plain style it inherited from A28. Apply the alignment
#    to a single cell first (keeps the style pool tidy - setting two
#    alignment properties at once on a single cell interns cleanly), then
#    fan the resulting format out to the rest of the row with a
#    formats-only paste so every cell in the row ends up sharing one style.
# ---------------------------------------------------------------------------
$b28 = $ws.Range("B28")
$b28.HorizontalAlignment = -4108   # xlCenter
$b28.VerticalAlignment = -4108     # xlCenter
$b28.Copy()
$ws.Range("C28:D28").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3) Shrink the two big merges so column E is no longer part of them
#    (B21:E27 -> B21:D27, B29:E35 -> B29:D35). Re-merging a range recomputes
#    per-cell borders/alignment for the whole rectangle, so stash the
#    original uniform formatting first and paint it back over the full
#    original footprint (including column E, which stays a plain cell)
#    once the merge geometry is updated.
# ---------------------------------------------------------------------------
$ws.Range("B21").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("B21:E27").UnMerge()
$ws.Range("B21:D27").Merge()
$ws.Range("Z100").Copy()
$ws.Range("B21:E27").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

$ws.Range("B29").Copy()
$ws.Range("Z100").PasteSpecial(-4122)
$ws.Range("B29:E35").UnMerge()
$ws.Range("B29:D35").Merge()
$ws.Range("Z100").Copy()
$ws.Range("B29:E35").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Move the active selection to D9 (cosmetic, matches the saved view state)
# ---------------------------------------------------------------------------
[void]$ws.Range("D9").Select()
